$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")
$ws.Rows.Item(20).Delete()
